$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VerifyCSVForExistingVersion")

# Row 6 (TestCaseID row for "server stand alone" case) now records a Fail
$ws.Range("H6").Value = "Fail"

# Remaining result cells (rows 8 through 44) are cleared out entirely
$ws.Range("H8:H44").ClearContents()
